# SectorGroup.xlsx - codeforIATI codelists deploy
#
# The columns D:G on Sheet1 hold, per row:
#   D = codeforiati:category-code
#   E = codeforiati:group-code
#   F = codeforiati:group-name
#   G = codeforiati:category-name
#
# The upstream codelist generator changed the column order for this sheet to:
#   D = codeforiati:group-name
#   E = codeforiati:category-name
#   F = codeforiati:group-code
#   G = codeforiati:category-code
#
# i.e. for every row (including the header row):
#   new D <- old F
#   new E <- old G
#   new F <- old E
#   new G <- old D
#
# Values must stay text (they are codelist codes like "110", "111", etc., not
# numbers), so the swap is done with Copy/PasteSpecial(xlPasteValues) via
# scratch columns rather than a plain Value assignment (which would let
# numeric-looking text turn into real numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$xlPasteValues = -4163

# 1) Stash the current D,E,F,G columns into scratch columns (Z,AA,AB,AC)
#    so the originals survive while we rebuild them in the new order.
$ws.Range("D1:D$lastRow").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteValues)

$ws.Range("E1:E$lastRow").Copy()
$ws.Range("AA1").PasteSpecial($xlPasteValues)

$ws.Range("F1:F$lastRow").Copy()
$ws.Range("AB1").PasteSpecial($xlPasteValues)

$ws.Range("G1:G$lastRow").Copy()
$ws.Range("AC1").PasteSpecial($xlPasteValues)

# 2) Paste them back in the new column order.
#    new D = old F (AB), new E = old G (AC), new F = old E (AA), new G = old D (Z)
$ws.Range("AB1:AB$lastRow").Copy()
$ws.Range("D1").PasteSpecial($xlPasteValues)

$ws.Range("AC1:AC$lastRow").Copy()
$ws.Range("E1").PasteSpecial($xlPasteValues)

$ws.Range("AA1:AA$lastRow").Copy()
$ws.Range("F1").PasteSpecial($xlPasteValues)

$ws.Range("Z1:Z$lastRow").Copy()
$ws.Range("G1").PasteSpecial($xlPasteValues)

# 3) Clean up the scratch columns.
$ws.Range("Z1:AC$lastRow").Clear()

$excel.CutCopyMode = $false
